$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date in column C for every existing data row (2-457)
for ($r = 2; $r -le 457; $r++) {
    $ws.Cells.Item($r, 3).Value = 45203
}

# Row 457 gains an explicit row height (matches rows above it)
$ws.Rows.Item(457).RowHeight = 15

# Append the new record as row 458
$ws.Cells.Item(458, 1).Value = "A 47146-2023"
$ws.Cells.Item(458, 2).Value = 45201
$ws.Cells.Item(458, 3).Value = 45203
$ws.Cells.Item(458, 4).Value = "VÄSTERBOTTENS LÄN"
$ws.Cells.Item(458, 5).Value = "ÅSELE"
$ws.Cells.Item(458, 7).Value = 0.8
$ws.Cells.Item(458, 8).Value = 0
$ws.Cells.Item(458, 9).Value = 0
$ws.Cells.Item(458, 10).Value = 0
$ws.Cells.Item(458, 11).Value = 0
$ws.Cells.Item(458, 12).Value = 0
$ws.Cells.Item(458, 13).Value = 0
$ws.Cells.Item(458, 14).Value = 0
$ws.Cells.Item(458, 15).Value = 0
$ws.Cells.Item(458, 16).Value = 0
$ws.Cells.Item(458, 17).Value = 0

# Match number formatting used by the date columns elsewhere in the sheet
$ws.Range("B458:C458").NumberFormat = "YYYY-MM-DD"

# The (empty) species column uses a wrap-text style throughout the sheet
$ws.Range("R458").WrapText = $true
